$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A - shifts all existing data right by one column
$ws.Columns("A:A").Insert()

# New header / value for the inserted "Code Article" column
$ws.Range("A2").Value = "Code Article"
$ws.Range("A3").Value = "E-COM11"

$wb.Save()
